$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the column headers for "Energy" (column B) and "NEIG" (column D).
$ws.Range("B1").Value = "NEIG"
$ws.Range("D1").Value = "Energy"

# Add new data row 55 (June 1, 2024) after the existing last row (54).
# Copy the date cell's format down from the row above so the new row matches
# the existing date styling instead of creating a brand-new style entry.
$ws.Range("A54").Copy()
$ws.Range("A55").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A55").Value = 45444
$ws.Range("B55").Value = -0.455
$ws.Range("C55").Value = 0.438
$ws.Range("D55").Value = -0.574
$ws.Range("E55").Value = 0.432
$ws.Range("F55").Value = 1.655
